$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is a brand-new appended row. Plain-looking text is written straight
# through; numeric-looking / empty strings need a quick NumberFormat or
# quote-prefix nudge so Excel's "user input" parser keeps them as text
# (matching the source rows above them), then the style is reset back to
# Normal so no stray cell format lingers on the new row.

$ws.Range("A8").Value = "'"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = "أحمد شريم"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "2323"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = "ايتا"
$ws.Range("E8").Value = "الرحلة 2"
$ws.Range("F8").Value = "C2"
$ws.Range("G8").Value = "NRC"
$ws.Range("H8").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٢٠:٥٧ م"
